$d = $word.ActiveDocument

$d.Content.Find.Execute("824÷9=91, 5", $true, $false, $false, $false, $false, $true, 1, $false, "266÷7=38, 0", 2) | Out-Null
$d.Content.Find.Execute("300÷7=42, 6", $true, $false, $false, $false, $false, $true, 1, $false, "847÷8=105, 7", 2) | Out-Null
$d.Content.Find.Execute("685÷3=228, 1", $true, $false, $false, $false, $false, $true, 1, $false, "625÷5=125, 0", 2) | Out-Null
$d.Content.Find.Execute("389÷8=48, 5", $true, $false, $false, $false, $false, $true, 1, $false, "434÷8=54, 2", 2) | Out-Null
$d.Content.Find.Execute("612÷3=204, 0", $true, $false, $false, $false, $false, $true, 1, $false, "582÷7=83, 1", 2) | Out-Null
$d.Content.Find.Execute("487÷8=60, 7", $true, $false, $false, $false, $false, $true, 1, $false, "485÷5=97, 0", 2) | Out-Null
$d.Content.Find.Execute("940÷3=313, 1", $true, $false, $false, $false, $false, $true, 1, $false, "847÷4=211, 3", 2) | Out-Null
$d.Content.Find.Execute("294÷6=49, 0", $true, $false, $false, $false, $false, $true, 1, $false, "317÷8=39, 5", 2) | Out-Null
$d.Content.Find.Execute("965÷5=193, 0", $true, $false, $false, $false, $false, $true, 1, $false, "349÷5=69, 4", 2) | Out-Null
$d.Content.Find.Execute("169÷3=56, 1", $true, $false, $false, $false, $false, $true, 1, $false, "247÷4=61, 3", 2) | Out-Null
$d.Content.Find.Execute("777÷8=97, 1", $true, $false, $false, $false, $false, $true, 1, $false, "436÷8=54, 4", 2) | Out-Null
$d.Content.Find.Execute("895÷5=179, 0", $true, $false, $false, $false, $false, $true, 1, $false, "267÷8=33, 3", 2) | Out-Null
$d.Content.Find.Execute("606÷7=86, 4", $true, $false, $false, $false, $false, $true, 1, $false, "461÷7=65, 6", 2) | Out-Null
$d.Content.Find.Execute("774÷3=258, 0", $true, $false, $false, $false, $false, $true, 1, $false, "264÷2=132, 0", 2) | Out-Null
$d.Content.Find.Execute("676÷7=96, 4", $true, $false, $false, $false, $false, $true, 1, $false, "469÷6=78, 1", 2) | Out-Null
$d.Content.Find.Execute("385÷4=96, 1", $true, $false, $false, $false, $false, $true, 1, $false, "649÷2=324, 1", 2) | Out-Null
$d.Content.Find.Execute("339÷3=113, 0", $true, $false, $false, $false, $false, $true, 1, $false, "193÷6=32, 1", 2) | Out-Null
$d.Content.Find.Execute("985÷6=164, 1", $true, $false, $false, $false, $false, $true, 1, $false, "691÷6=115, 1", 2) | Out-Null
$d.Content.Find.Execute("570÷4=142, 2", $true, $false, $false, $false, $false, $true, 1, $false, "483÷6=80, 3", 2) | Out-Null
$d.Content.Find.Execute("140÷4=35, 0", $true, $false, $false, $false, $false, $true, 1, $false, "426÷5=85, 1", 2) | Out-Null
$d.Content.Find.Execute("896÷6=149, 2", $true, $false, $false, $false, $false, $true, 1, $false, "860÷9=95, 5", 2) | Out-Null
$d.Content.Find.Execute("327÷8=40, 7", $true, $false, $false, $false, $false, $true, 1, $false, "219÷2=109, 1", 2) | Out-Null
$d.Content.Find.Execute("377÷5=75, 2", $true, $false, $false, $false, $false, $true, 1, $false, "262÷6=43, 4", 2) | Out-Null
$d.Content.Find.Execute("329÷8=41, 1", $true, $false, $false, $false, $false, $true, 1, $false, "658÷9=73, 1", 2) | Out-Null
$d.Content.Find.Execute("352÷6=58, 4", $true, $false, $false, $false, $false, $true, 1, $false, "877÷5=175, 2", 2) | Out-Null
